# Fix a typo in the sheet name "2-3" -> "2-2", and make that sheet the
# active/selected tab (it was previously "team" that held the selection).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2-3")
$ws.Name = "2-2"

# Moves the active/selected tab from "team" to this (renamed) sheet,
# matching the workbook's activeTab/tabSelected bookkeeping.
$ws.Activate()
